# TC01_Canine_Filter_NeutStatus-Yes.xlsx
# Fix the "CasesTab" query (startup!B2): it incorrectly returned an extra
# `Cohort` column that isn't produced anywhere else in this workbook's
# other tab queries. Drop the trailing `co.cohort_description` projection
# (and the now-dangling trailing comma on the prior line) so the query
# matches the corrected report definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.neutered_indicator IN ["Yes"]  
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCasesQuery

# Reflect the author's final selection/scroll position on the sheet.
$ws.Range("B2").Select()
